$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Reports icon SVG markup reused for the new Special Reports row ---
$reportsIcon = @'
<g transform='translate(-2.2,-2.1) scale(1.18)'>
                                <path opacity='0.5' d='M21 15.9983V9.99826C21 7.16983 21 5.75562 20.1213 4.87694C19.3529 4.10856 18.175 4.01211 16 4H8C5.82497 4.01211 4.64706 4.10856 3.87868 4.87694C3 5.75562 3 7.16983 3 9.99826V15.9983C3 18.8267 3 20.2409 3.87868 21.1196C4.75736 21.9983 6.17157 21.9983 9 21.9983H15C17.8284 21.9983 19.2426 21.9983 20.1213 21.1196C21 20.2409 21 18.8267 21 15.9983Z' fill='currentColor'></path>
                                <path d='M8 3.5C8 2.67157 8.67157 2 9.5 2H14.5C15.3284 2 16 2.67157 16 3.5V4.5C16 5.32843 15.3284 6 14.5 6H9.5C8.67157 6 8 5.32843 8 4.5V3.5Z' fill='currentColor'></path>
                                <path fill-rule='evenodd' clip-rule='evenodd' d='M6.25 10.5C6.25 10.0858 6.58579 9.75 7 9.75H17C17.4142 9.75 17.75 10.0858 17.75 10.5C17.75 10.9142 17.4142 11.25 17 11.25H7C6.58579 11.25 6.25 10.9142 6.25 10.5ZM7.25 14C7.25 13.5858 7.58579 13.25 8 13.25H16C16.4142 13.25 16.75 13.5858 16.75 14C16.75 14.4142 16.4142 14.75 16 14.75H8C7.58579 14.75 7.25 14.4142 7.25 14ZM8.25 17.5C8.25 17.0858 8.58579 16.75 9 16.75H15C15.4142 16.75 15.75 17.0858 15.75 17.5C15.75 17.9142 15.4142 18.25 15 18.25H9C8.58579 18.25 8.25 17.9142 8.25 17.5Z' fill='currentColor'></path>
                                </g>
'@

# --- New row 28: Special Reports module ---
$ws.Range("C28").Value = "Laporan Khas"
$ws.Range("D28").Value = "Reports Special"
$ws.Range("E28").Value = "Generates special report and displays data summaries on system usage and performance."
$ws.Range("A28").Value = "special_reports"
$ws.Range("F28").Value = $reportsIcon
$ws.Range("G28").Value = 6

# The multi-line icon markup makes Excel flag the row with a custom height;
# re-running AutoFit settles it back to the sheet's standard row height.
$ws.Rows.Item(28).AutoFit()

# --- Update the "API Logs" row (row 26) to "Backend Logs" ---
$ws.Range("C26").Value = "Backend Logs"
$ws.Range("D26").Value = "Backend Logs"

# --- Row 23 (Knowledge Base) order number bumped from 6 to 7 ---
$ws.Range("G23").Value = 7

# --- Match the saved selection state from the edit ---
$ws.Range("D24").Select()
